# Review pesticide names (close #333)
# Rename the abbreviated pesticide/chemical names on the "Chemical Properties"
# sheet (columns A and W) to their full names. New shared strings are
# appended automatically by the engine when the .Value is set to text that
# doesn't already exist in the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chemical Properties")

# Row 3: Azoxy -> Azoxystrobin
$ws.Range("A3").Value = "Azoxystrobin"
$ws.Range("W3").Value = "Azoxystrobin"

# Row 4: Benta -> Bentazone
$ws.Range("A4").Value = "Bentazone"
$ws.Range("W4").Value = "Bentazone"

# Row 5: Cicloxidim -> Cycloxydim
$ws.Range("A5").Value = "Cycloxydim"
$ws.Range("W5").Value = "Cycloxydim"

# Row 6: Cyhalo -> Cyhalofop-butyl (column A); column W keeps mirroring the
# (pre-existing, unchanged) value from row 5 -> Cycloxydim
$ws.Range("A6").Value = "Cyhalofop-butyl"
$ws.Range("W6").Value = "Cycloxydim"

# Row 7: Difeno -> Difenoconazole
$ws.Range("A7").Value = "Difenoconazole"
$ws.Range("W7").Value = "Difenoconazole"

# Row 8 (MCPA) is left untouched.

# Row 9: Penoxulam -> Penoxsulam
$ws.Range("A9").Value = "Penoxsulam"
$ws.Range("W9").Value = "Penoxsulam"

# Cosmetic view-state update captured in the diff: the sheet was re-zoomed.
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 175
